$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-unused columns H:L entirely (was episodes 6-10)
$ws.Range("H1:L2").Clear()

# Row 1 header changes
$ws.Range("B1").Value = "Values"

# Row 2: change B2 from number to text label, and update C2:G2 values
$ws.Range("B2").Value = "Final Value"
$ws.Range("C2").Value = 2651718.026057291
$ws.Range("D2").Value = 2861276.966189215
$ws.Range("E2").Value = 1858585.627848443
$ws.Range("F2").Value = 2665675.617408922
$ws.Range("G2").Value = 2403791.3239417

# New row 3
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "Annualized Return"
$ws.Range("C3").Value = 0.3800047767090269
$ws.Range("D3").Value = 0.4151107341519846
$ws.Range("E3").Value = 0.2271688888706627
$ws.Range("F3").Value = 0.3823996129066523
$ws.Range("G3").Value = 0.3359824892542176

# New row 4
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "Sharpe Ratio"
$ws.Range("C4").Value = 0.7277318355804748
$ws.Range("D4").Value = 1.074122742102152
$ws.Range("E4").Value = 0.5918018506522419
$ws.Range("F4").Value = 0.7356774356629673
$ws.Range("G4").Value = 0.6779838054083842

# Apply the same style (s=1, centered bold border) as A2 to A3 and A4
$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("A4").PasteSpecial(-4122) # xlPasteFormats
